# "UP to date gemaakt"
# Add a new logboek week sheet ("week 4") right before the "Totaal" sheet.
# It is a copy of the "week 2" sheet layout with its first entry's data
# updated, and rows for the entries that did not happen that week cleared.

$wb = $excel.ActiveWorkbook

$weekTwo = $wb.Worksheets.Item("week 2")

# Copy "week 2" and place the copy directly after it -> ends up right
# before "Totaal", becomes the active sheet/tab automatically (matching
# the tabSelected move and activeTab=4 change in the workbook).
$weekTwo.Copy($null, $weekTwo) | Out-Null

$weekFour = $wb.Worksheets.Item("week 2 (2)")
$weekFour.Name = "week 4"

# Update the first logboek entry for the new week.
$weekFour.Range("B7").Value2 = 41669
$weekFour.Range("C7").Value2 = 0.40972222222222227
$weekFour.Range("D7").Value2 = 0.42708333333333331
$weekFour.Range("F7").Value2 = "Bijgewerkt naar leraars project "

# The remaining entries from the copied "week 2" sheet don't apply to
# this week, so clear them out.
$weekFour.Range("C8:F8").ClearContents() | Out-Null
$weekFour.Range("C9:F9").ClearContents() | Out-Null
$weekFour.Range("A11:F11").ClearContents() | Out-Null
$weekFour.Range("C12:F12").ClearContents() | Out-Null
$weekFour.Range("C13:F13").ClearContents() | Out-Null

# Those rows had extra height for wrapped text; restore default height
# now that they are empty again.
$weekFour.Rows.Item(8).AutoFit() | Out-Null
$weekFour.Rows.Item(9).AutoFit() | Out-Null
$weekFour.Rows.Item(11).AutoFit() | Out-Null
$weekFour.Rows.Item(13).AutoFit() | Out-Null

# Make the new sheet active with its selection on E9.
$weekFour.Activate() | Out-Null
$weekFour.Range("E9").Select() | Out-Null
